# ---------------------------------------------------------------------------
# Edit script: apply the changes described by the days.xlsx diff.
#  1. Add 3 new shared strings (done implicitly by writing the new text values)
#  2. AP228: 7 -> 8
#  3. Row 236: height 12.8 -> 12.85 ; clear X236 (formula removed)
#  4. Append new rows 237, 238, 239, 240 with their data / formulas
#  5. Update the sheet view (zoom, selection) as far as the object model allows
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- small helpers to reproduce the workbook's existing cell styles -------
# style "2"  -> integer custom format   [$-409]0
# style "3"  -> 2 decimal custom format [$-409]0.00
# style "4"  -> 2 decimal format        #,##0.00
# style "7"  -> General + wrap text
# style "13" -> General, Calibri 11 font (black)
function Set-IntStyle($rng) {
    $rng.NumberFormat = "[$-409]0"
}
function Set-Float2Style($rng) {
    $rng.NumberFormat = "[$-409]0.00"
}
function Set-ThousandsStyle($rng) {
    $rng.NumberFormat = "#,##0.00"
}
function Set-WrapStyle($rng) {
    $rng.WrapText = $true
}
function Set-CalibriStyle($rng) {
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Color = 0
}

# =====================================================================
# 1. AP228 7 -> 8
# =====================================================================
$ws.Range("AP228").Value = 8

# =====================================================================
# 2. Row 236: new row height, X236 formula removed
# =====================================================================
$ws.Rows.Item(236).RowHeight = 12.85
$ws.Range("X236").ClearContents()
$ws.Range("X236").NumberFormat = "General"

# =====================================================================
# 3. Row 237 (new)
# =====================================================================
$ws.Range("A237").Formula = "=A236+1"
$ws.Range("B237").Value = 44115.4951388889
$ws.Range("C237").Value = 0
$ws.Range("D237").Value = "Weather forecast"
$ws.Range("F237").Value = "Fair"
$ws.Range("G237").Value = 71
Set-IntStyle $ws.Range("G237")
$ws.Range("H237").Value = 61
Set-IntStyle $ws.Range("H237")
$ws.Range("I237").Value = 72
Set-IntStyle $ws.Range("I237")
$ws.Range("J237").Value = "NNW "
Set-IntStyle $ws.Range("J237")
$ws.Range("K237").Value = 9
Set-IntStyle $ws.Range("K237")
$ws.Range("L237").Value = 0
Set-IntStyle $ws.Range("L237")
$ws.Range("M237").Value = "Good"
$ws.Range("N237").Value = 0

# =====================================================================
# 4. Row 238 (new)
# =====================================================================
$ws.Rows.Item(238).RowHeight = 13.8
$ws.Range("A238").Formula = "=A237+1"
$ws.Range("B238").Value = 44116.4597222222
$ws.Range("C238").Value = 1
$ws.Range("F238").Value = "Fair"
$ws.Range("G238").Formula = "=(75+79)/2"
Set-IntStyle $ws.Range("G238")
$ws.Range("H238").Value = 42
Set-IntStyle $ws.Range("H238")
$ws.Range("I238").Value = 28
Set-IntStyle $ws.Range("I238")
$ws.Range("J238").Value = "N"
Set-IntStyle $ws.Range("J238")
$ws.Range("K238").Value = 22
Set-IntStyle $ws.Range("K238")
$ws.Range("L238").Value = 30
Set-IntStyle $ws.Range("L238")
$ws.Range("M238").Value = "Good"
$ws.Range("N238").Value = 0
$ws.Range("O238").Value = "Sketchers 2"
$ws.Range("P238").Value = "Green Oaks - Little Road"
Set-CalibriStyle $ws.Range("P238")
$ws.Range("Q238").Value = 6.85
Set-Float2Style $ws.Range("Q238")
$ws.Range("R238").Value = 860
Set-IntStyle $ws.Range("R238")
$ws.Range("S238").Formula = "=T238+R238"
Set-IntStyle $ws.Range("S238")
$ws.Range("T238").Value = 14920
Set-IntStyle $ws.Range("T238")
$ws.Range("U238").Formula = "=(120+3)/60"
Set-Float2Style $ws.Range("U238")
$ws.Range("V238").Formula = "=(120+16)/60"
Set-Float2Style $ws.Range("V238")
$ws.Range("W238").Formula = "=V238-U238"
Set-Float2Style $ws.Range("W238")
$ws.Range("X238").Formula = "=Q236/U236"
Set-Float2Style $ws.Range("X238")
$ws.Range("Y238").Value = 1
$ws.Range("Z238").Formula = "=Q238/Y238"
Set-Float2Style $ws.Range("Z238")
$ws.Range("AA238").Formula = "=17+56/60"
Set-Float2Style $ws.Range("AA238")
$ws.Range("AB238").Value = 292
Set-IntStyle $ws.Range("AB238")
$ws.Range("AC238").Value = 555
Set-IntStyle $ws.Range("AC238")
$ws.Range("AD238").Value = 69
Set-IntStyle $ws.Range("AD238")
$ws.Range("AE238").Value = 110
$ws.Range("AF238").Formula = "=16+52/60"
Set-Float2Style $ws.Range("AF238")
$ws.Range("AG238").Formula = "=17+31/60"
Set-Float2Style $ws.Range("AG238")
$ws.Range("AH238").Formula = "=17+54/60"
Set-Float2Style $ws.Range("AH238")
$ws.Range("AI238").Formula = "=18+51/60"
Set-Float2Style $ws.Range("AI238")
$ws.Range("AJ238").Formula = "=18+14/60"
Set-Float2Style $ws.Range("AJ238")
$ws.Range("AK238").Formula = "=17+54/60"
Set-Float2Style $ws.Range("AK238")
$ws.Range("AL238").Formula = "=60/3.3"
Set-Float2Style $ws.Range("AL238")
$ws.Range("AM238").NumberFormat = "General"
$ws.Range("AP238").Value = 8
Set-IntStyle $ws.Range("AP238")
$ws.Range("AQ238").Value = 2
$ws.Range("AR238").Value = 0
$ws.Range("AS238").Value = 0
$ws.Range("AT238").Formula = "=60*U238-SUM(AU238:AY238)"
Set-ThousandsStyle $ws.Range("AT238")
$ws.Range("AU238").Formula = "=6+18/60"
Set-Float2Style $ws.Range("AU238")
$ws.Range("AV238").Formula = "=16/60"
Set-Float2Style $ws.Range("AV238")
$ws.Range("AW238").Value = 0
Set-Float2Style $ws.Range("AW238")
$ws.Range("AX238").Value = 0
Set-Float2Style $ws.Range("AX238")
$ws.Range("AY238").Value = 0
Set-Float2Style $ws.Range("AY238")
$ws.Range("AZ238").Value = "Garman vivoactive 3"
$ws.Range("BA238").Value = "Connect"
$ws.Range("BB238").Value = 0

# =====================================================================
# 5. Row 239 (new)
# =====================================================================
$ws.Rows.Item(239).RowHeight = 12.85
$ws.Range("A239").Value = 771
$ws.Range("B239").Value = 44117.5930555556
$ws.Range("C239").Value = 1
$ws.Range("F239").Value = "Partly Cloudy "
Set-WrapStyle $ws.Range("F239")
$ws.Range("G239").Formula = "=(77+81)/2"
Set-IntStyle $ws.Range("G239")
$ws.Range("H239").Value = 42
Set-IntStyle $ws.Range("H239")
$ws.Range("I239").Value = 26
Set-IntStyle $ws.Range("I239")
$ws.Range("J239").Value = "SSW"
Set-IntStyle $ws.Range("J239")
$ws.Range("K239").Value = 6
Set-IntStyle $ws.Range("K239")
$ws.Range("L239").Value = 0
Set-IntStyle $ws.Range("L239")
$ws.Range("M239").Value = "Good"
$ws.Range("N239").Value = 0
$ws.Range("O239").Value = "Sketchers 2"
$ws.Range("P239").Value = "North Green Oaks"
$ws.Range("Q239").Value = 7.71
Set-Float2Style $ws.Range("Q239")
$ws.Range("R239").Value = 500
Set-IntStyle $ws.Range("R239")
$ws.Range("S239").Value = 17039
Set-IntStyle $ws.Range("S239")
$ws.Range("T239").Formula = "=S239-R239"
Set-IntStyle $ws.Range("T239")
$ws.Range("U239").Formula = "=(120+6)/60"
Set-Float2Style $ws.Range("U239")
$ws.Range("V239").Formula = "=(120+15)/60"
Set-Float2Style $ws.Range("V239")
$ws.Range("W239").Formula = "=V239-U239"
Set-Float2Style $ws.Range("W239")
$ws.Range("X239").Formula = "=Q239/U239"
Set-Float2Style $ws.Range("X239")
$ws.Range("Y239").Value = 1
$ws.Range("Z239").Formula = "=Q239/Y239"
Set-Float2Style $ws.Range("Z239")
$ws.Range("AA239").Formula = "=16+24/60"
Set-Float2Style $ws.Range("AA239")
$ws.Range("AB239").Value = 1037
Set-IntStyle $ws.Range("AB239")
$ws.Range("AC239").Value = 820
Set-IntStyle $ws.Range("AC239")
$ws.Range("AD239").Value = 128
Set-IntStyle $ws.Range("AD239")
$ws.Range("AE239").Value = 141
$ws.Range("AF239").Formula = "=16+31/60"
Set-Float2Style $ws.Range("AF239")
$ws.Range("AG239").Formula = "=15+58/60"
Set-Float2Style $ws.Range("AG239")
$ws.Range("AH239").Formula = "=16+7/60"
Set-Float2Style $ws.Range("AH239")
$ws.Range("AI239").Formula = "=16+38/60"
Set-Float2Style $ws.Range("AI239")
$ws.Range("AJ239").Formula = "=16+26/60"
Set-Float2Style $ws.Range("AJ239")
$ws.Range("AK239").Formula = "=15+54/60"
Set-Float2Style $ws.Range("AK239")
$ws.Range("AL239").Formula = "=60/3.7"
Set-Float2Style $ws.Range("AL239")
$ws.Range("AP239").Value = 11
Set-IntStyle $ws.Range("AP239")
$ws.Range("AQ239").Value = 2
$ws.Range("AR239").Value = 0
$ws.Range("AS239").Value = 0
$ws.Range("AT239").Formula = "=60*U239-SUM(AU239:AY239)"
Set-ThousandsStyle $ws.Range("AT239")
$ws.Range("AU239").Formula = "=4+41/60"
Set-Float2Style $ws.Range("AU239")
$ws.Range("AV239").Formula = "=26+32/60"
Set-Float2Style $ws.Range("AV239")
$ws.Range("AW239").Formula = "=(95+14/60)"
Set-Float2Style $ws.Range("AW239")
$ws.Range("AX239").Value = 0
Set-Float2Style $ws.Range("AX239")
$ws.Range("AY239").Value = 0
Set-Float2Style $ws.Range("AY239")
$ws.Range("AZ239").Value = "Garman vivoactive 3"
$ws.Range("BA239").Value = "Connect"
$ws.Range("BB239").Value = 0

# =====================================================================
# 6. Row 240 (new) -- no A240 value (id left blank)
# =====================================================================
$ws.Range("B240").Value = 44118.4534722222
$ws.Range("C240").Value = 1
$ws.Range("F240").Value = "Fair"
$ws.Range("G240").Formula = "=79+11/60*(83-79)"
Set-IntStyle $ws.Range("G240")
$ws.Range("H240").Formula = "=60+11/60*(64-60)"
Set-IntStyle $ws.Range("H240")
$ws.Range("I240").Formula = "=56+11/60+(60-63)"
Set-IntStyle $ws.Range("I240")
$ws.Range("J240").Value = "S"
Set-IntStyle $ws.Range("J240")
$ws.Range("K240").Value = 13
Set-IntStyle $ws.Range("K240")
$ws.Range("L240").Value = 28
Set-IntStyle $ws.Range("L240")
$ws.Range("M240").Value = "Good"
$ws.Range("N240").Value = 1
$ws.Range("O240").Value = "Sketchers 2"
$ws.Range("P240").Value = "Stoval Park"
$ws.Range("Q240").Value = 4.2
Set-Float2Style $ws.Range("Q240")
$ws.Range("T240").Value = 9478
Set-IntStyle $ws.Range("T240")
$ws.Range("U240").Value = 1.11
Set-Float2Style $ws.Range("U240")
$ws.Range("X240").Formula = "=Q240/U240"
Set-Float2Style $ws.Range("X240")
$ws.Range("Y240").Value = 4
$ws.Range("Z240").Formula = "=Q240/Y240"
Set-Float2Style $ws.Range("Z240")
$ws.Range("AA240").Value = 15.9
Set-Float2Style $ws.Range("AA240")
$ws.Range("AB240").Value = 98.5
Set-IntStyle $ws.Range("AB240")
$ws.Range("AC240").Value = 552
Set-IntStyle $ws.Range("AC240")
$ws.Range("AP240").Value = 0
Set-IntStyle $ws.Range("AP240")
$ws.Range("AQ240").Value = 1
$ws.Range("AR240").Value = 1
$ws.Range("AS240").Value = 0
$ws.Range("AZ240").Value = "Garman vivoactive 3"
$ws.Range("BA240").Value = "Connect"
$ws.Range("BB240").Value = 1
$ws.Range("BC240").Value = "Farmin not charged"

# =====================================================================
# 7. Sheet view: zoom + final selection (best effort - engine only
#    models a single active selection / pane topLeftCell)
# =====================================================================
$win = $excel.ActiveWindow
$win.Zoom = 110
$ws.Range("BC241").Select()
